$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, $Addr, $Text)
    $cell = $Sheet.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-CellText $ws 'D2' '69.113.81'
Set-CellText $ws 'E2' '  +0.25%  '
Set-CellText $ws 'D3' '3.746.75'
Set-CellText $ws 'E3' '  +0.17%  '
Set-CellText $ws 'E4' '  +0.04%  '
Set-CellText $ws 'D5' '601.53'
Set-CellText $ws 'E5' '  +0.01%  '
Set-CellText $ws 'D6' '167.09'
Set-CellText $ws 'E6' '  -0.51%  '
Set-CellText $ws 'D7' '3.747.87'
Set-CellText $ws 'E7' '  +0.24%  '
Set-CellText $ws 'E8' '  -0.02%  '
Set-CellText $ws 'E9' '  +1.05%  '
Set-CellText $ws 'D10' '0.169'
Set-CellText $ws 'E10' '  +2.43%  '
Set-CellText $ws 'E11' '  +1.02%  '
Set-CellText $ws 'E12' '  +0.03%  '
Set-CellText $ws 'D13' '37.95'
Set-CellText $ws 'E13' '  -0.65%  '
Set-CellText $ws 'E14' '  +1.50%  '
Set-CellText $ws 'D15' '4.377.00'
Set-CellText $ws 'E15' '  +0.19%  '
Set-CellText $ws 'D16' '3.745.37'
Set-CellText $ws 'E16' '  -0.04%  '
Set-CellText $ws 'D17' '69.130.62'
Set-CellText $ws 'E17' '  +0.34%  '
Set-CellText $ws 'E18' '  +1.29%  '
Set-CellText $ws 'D19' '17.39'
Set-CellText $ws 'E19' '  +0.88%  '
Set-CellText $ws 'D20' '0.113'
Set-CellText $ws 'E20' '  -1.61%  '
Set-CellText $ws 'D21' '11.05'
Set-CellText $ws 'E21' '  +7.81%  '
Set-CellText $ws 'D22' '492.94'
Set-CellText $ws 'E22' '  -0.75%  '
Set-CellText $ws 'E23' '  +0.61%  '
Set-CellText $ws 'E24' '  +7.30%  '
Set-CellText $ws 'D25' '84.87'
Set-CellText $ws 'E25' '  -0.43%  '
Set-CellText $ws 'E26' '  -0.25%  '
Set-CellText $ws 'D27' '12.27'
Set-CellText $ws 'E27' '  -0.04%  '
Set-CellText $ws 'D28' '10.07'
Set-CellText $ws 'E28' '  -0.52%  '
Set-CellText $ws 'E30' '  +1.01%  '
Set-CellText $ws 'B31' 'NEARProtocol'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-CellText $ws 'D31' '8.13'
Set-CellText $ws 'E31' '  +2.59%  '
Set-CellText $ws 'B32' 'ImmutableX'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 'D32' '2.47'
Set-CellText $ws 'E32' '  +2.03%  '
Set-CellText $ws 'D33' '31.49'
Set-CellText $ws 'E33' '  -0.79%  '
Set-CellText $ws 'D34' '3.897.89'
Set-CellText $ws 'E34' '  +0.27%  '
Set-CellText $ws 'D35' '3.683.66'
Set-CellText $ws 'E35' '  +0.30%  '
Set-CellText $ws 'E36' '  -0.13%  '
Set-CellText $ws 'E37' '  +0.08%  '
Set-CellText $ws 'D38' '5.96'
Set-CellText $ws 'E38' '  +2.71%  '
Set-CellText $ws 'E39' '  -0.10%  '
Set-CellText $ws 'D40' '0.137'
Set-CellText $ws 'E40' '  +2.53%  '
Set-CellText $ws 'E41' '  +0.50%  '
Set-CellText $ws 'E42' '  +5.32%  '
Set-CellText $ws 'D43' '48.64'
Set-CellText $ws 'E43' '  -0.77%  '
Set-CellText $ws 'D44' '425.67'
Set-CellText $ws 'E44' '  -2.76%  '
Set-CellText $ws 'E45' '  -0.68%  '
Set-CellText $ws 'E46' '  +0.64%  '
Set-CellText $ws 'D48' '40.16'
Set-CellText $ws 'E48' '  -0.75%  '
Set-CellText $ws 'D49' '141.48'
Set-CellText $ws 'E49' '  -0.73%  '
Set-CellText $ws 'D50' '2.791.24'
Set-CellText $ws 'E50' '  +1.54%  '
Set-CellText $ws 'D51' '0.0353'
Set-CellText $ws 'E51' '  +0.35%  '
